# stats & report adjustment
# Merge the separate "credit" placeholders into combined "credit + code" placeholders
# on the "5、储值卡核销收入" row, and move the active selection to C11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "{{creditAndCodeAmount}}"
$ws.Range("C10").Value = "{{creditAndCodeAmountM}}"

$ws.Range("C11").Select()
